$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.5
$summary.Range("B4").Value = -2.49
$summary.Range("B5").Value = -0.92
$summary.Range("B6").Value = 54
$summary.Range("B7").Value = 22
$summary.Range("B9").Value = 40.74

# --- Sheet: Strategy Status (row 4 = MarketMaking) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.5
$status.Range("D4").Value = 54
$status.Range("E4").Value = -2.49
$status.Range("F4").Value = -2.5
$status.Range("G4").Value = 40.74

# --- Sheet: All Trades (append new trade row 55) ---
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(55, 1).Value = 54
# "2026-02-17" looks like a date, so Excel would auto-convert it to a date
# serial. Force it to stay literal text (matching the source data) without
# leaving a lingering text number-format on the cell.
$allTrades.Cells.Item(55, 2).NumberFormat = "@"
$allTrades.Cells.Item(55, 2).Value = "2026-02-17"
$allTrades.Cells.Item(55, 2).ClearFormats()
$allTrades.Cells.Item(55, 3).Value = "13:29:25"
$allTrades.Cells.Item(55, 4).Value = "MarketMaking"
$allTrades.Cells.Item(55, 5).Value = "DOWN"
$allTrades.Cells.Item(55, 6).Value = 0.02
$allTrades.Cells.Item(55, 7).Value = 0.03
$allTrades.Cells.Item(55, 8).Value = "CLOSED"
$allTrades.Cells.Item(55, 9).Value = 50
$allTrades.Cells.Item(55, 10).Value = 0.01
$allTrades.Cells.Item(55, 11).Value = 97.5
$allTrades.Cells.Item(55, 12).Value = 0
$allTrades.Cells.Item(55, 13).Value = 0
$allTrades.Cells.Item(55, 14).Value = 0.6
$allTrades.Cells.Item(55, 15).Value = "Normal spread capture: 19600 bps"
$allTrades.Cells.Item(55, 16).Value = "early_exit"
$allTrades.Cells.Item(55, 17).Value = 0.11

# --- Sheet: MarketMaking (append same new trade row 55) ---
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(55, 1).Value = 54
$mm.Cells.Item(55, 2).NumberFormat = "@"
$mm.Cells.Item(55, 2).Value = "2026-02-17"
$mm.Cells.Item(55, 2).ClearFormats()
$mm.Cells.Item(55, 3).Value = "13:29:25"
$mm.Cells.Item(55, 4).Value = "MarketMaking"
$mm.Cells.Item(55, 5).Value = "DOWN"
$mm.Cells.Item(55, 6).Value = 0.02
$mm.Cells.Item(55, 7).Value = 0.03
$mm.Cells.Item(55, 8).Value = "CLOSED"
$mm.Cells.Item(55, 9).Value = 50
$mm.Cells.Item(55, 10).Value = 0.01
$mm.Cells.Item(55, 11).Value = 97.5
$mm.Cells.Item(55, 12).Value = 0
$mm.Cells.Item(55, 13).Value = 0
$mm.Cells.Item(55, 14).Value = 0.6
$mm.Cells.Item(55, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(55, 16).Value = "early_exit"
$mm.Cells.Item(55, 17).Value = 0.11
